$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 55
$ws.Cells.Item($row, 1).Value = "Golang Tech Lead"
$ws.Cells.Item($row, 2).Value = "https://www.dice.com/job-detail/2eae25a3-93f1-44f6-a98a-5af568612ca2"
$ws.Cells.Item($row, 3).Value = "Tampa, Florida"
$ws.Cells.Item($row, 4).Value = "Third Party, Contract"
$ws.Cells.Item($row, 5).Value = "Depends on Experience"
$ws.Cells.Item($row, 6).Value = "Concent Software Solution LLC"
